# Rename the inline pictures' shape names in the document's
# headers/footers, per the commit:
#   footer1: image1.png -> image2.png  (wp:docPr/pic:cNvPr id="1"/"0")
#   footer2: image1.png -> image2.png  (wp:docPr/pic:cNvPr id="2"/"0")
#   header2: image2.jpg -> image1.jpg  (wp:docPr/pic:cNvPr id="3"/"0")
#
# InlineShape has no writable "Name" property in the Word object model,
# so the standard COM-interop idiom is to flip the picture to a floating
# Shape (which does expose .Name), rename it, then flip it back to an
# inline shape so the layout/anchoring is unchanged.

$d = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-InlinePicture($range, $newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

# Footer 1 (first footer story) - PearsonLogo, id 1
Rename-InlinePicture $sec.Footers.Item(1).Range "image2.png"

# Footer 2 (even-page footer story) - PearsonLogo, id 2
Rename-InlinePicture $sec.Footers.Item(2).Range "image2.png"

# Header 2 (even-page header story) - BTec_Logo-Orange, id 3
Rename-InlinePicture $sec.Headers.Item(2).Range "image1.jpg"
